$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# F1 header: copy formatting from E1 (reuses the same header style), then set text
$ws.Range("E1").Copy()
$ws.Range("F1").PasteSpecial(-4122)
$ws.Range("F1").Value = "time_taken"

# Fill time_taken values for data rows 2-86 as text (inline/shared strings, no special style)
$timeTaken = @(
    "2021-10-05 13:41:36.472961",
    "2021-10-05 13:41:36.472973",
    "2021-10-05 13:41:36.472977",
    "2021-10-05 13:41:36.472979",
    "2021-10-05 13:41:36.472982",
    "2021-10-05 13:41:36.472985",
    "2021-10-05 13:41:36.472987",
    "2021-10-05 13:41:36.472990",
    "2021-10-05 13:41:36.472992",
    "2021-10-05 13:41:36.472995",
    "2021-10-05 13:41:36.472998",
    "2021-10-05 13:41:36.473000",
    "2021-10-05 13:41:36.473002",
    "2021-10-05 13:41:36.473005",
    "2021-10-05 13:41:36.473007",
    "2021-10-05 13:41:36.473010",
    "2021-10-05 13:41:36.473012",
    "2021-10-05 13:41:36.473015",
    "2021-10-05 13:41:36.473017",
    "2021-10-05 13:41:36.473020",
    "2021-10-05 13:41:36.473022",
    "2021-10-05 13:41:36.473025",
    "2021-10-05 13:41:36.473027",
    "2021-10-05 13:41:36.473030",
    "2021-10-05 13:41:36.473032",
    "2021-10-05 13:41:36.473035",
    "2021-10-05 13:41:36.473037",
    "2021-10-05 13:41:36.473040",
    "2021-10-05 13:41:36.473043",
    "2021-10-05 13:41:36.473045",
    "2021-10-05 13:41:36.473048",
    "2021-10-05 13:41:36.473050",
    "2021-10-05 13:41:36.473053",
    "2021-10-05 13:41:36.473055",
    "2021-10-05 13:41:36.473058",
    "2021-10-05 13:41:36.473060",
    "2021-10-05 13:41:36.473063",
    "2021-10-05 13:41:36.473065",
    "2021-10-05 13:41:36.473068",
    "2021-10-05 13:41:36.473071",
    "2021-10-05 13:41:36.473073",
    "2021-10-05 13:41:36.473076",
    "2021-10-05 13:41:36.473078",
    "2021-10-05 13:41:36.473081",
    "2021-10-05 13:41:36.473083",
    "2021-10-05 13:41:36.473086",
    "2021-10-05 13:41:36.473088",
    "2021-10-05 13:41:36.473090",
    "2021-10-05 13:41:36.473093",
    "2021-10-05 13:41:36.473095",
    "2021-10-05 13:41:36.473098",
    "2021-10-05 13:41:36.473100",
    "2021-10-05 13:41:36.473103",
    "2021-10-05 13:41:36.473106",
    "2021-10-05 13:41:36.473108",
    "2021-10-05 13:41:36.473111",
    "2021-10-05 13:41:36.473113",
    "2021-10-05 13:41:36.473115",
    "2021-10-05 13:41:36.473118",
    "2021-10-05 13:41:36.473120",
    "2021-10-05 13:41:36.473123",
    "2021-10-05 13:41:36.473125",
    "2021-10-05 13:41:36.473128",
    "2021-10-05 13:41:36.473130",
    "2021-10-05 13:41:36.473133",
    "2021-10-05 13:41:36.473136",
    "2021-10-05 13:41:36.473139",
    "2021-10-05 13:41:36.473141",
    "2021-10-05 13:41:36.473144",
    "2021-10-05 13:41:36.473146",
    "2021-10-05 13:41:36.473149",
    "2021-10-05 13:41:36.473151",
    "2021-10-05 13:41:36.473154",
    "2021-10-05 13:41:36.473156",
    "2021-10-05 13:41:36.473158",
    "2021-10-05 13:41:36.473161",
    "2021-10-05 13:41:36.473165",
    "2021-10-05 13:41:36.473168",
    "2021-10-05 13:41:36.473171",
    "2021-10-05 13:41:36.473173",
    "2021-10-05 13:41:36.473176",
    "2021-10-05 13:41:36.473178",
    "2021-10-05 13:41:36.473181",
    "2021-10-05 13:41:36.473183",
    "2021-10-05 13:41:36.473186"
)

for ($i = 0; $i -lt $timeTaken.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 6).Value = $timeTaken[$i]
}
